$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "30.727.93"
$ws.Cells.Item(2, 5).Value = "  +0.76%  "
$ws.Cells.Item(3, 4).Value = "2.148.27"
$ws.Cells.Item(3, 5).Value = "  +1.92%  "
$ws.Cells.Item(4, 5).Value = "  +0.51%  "
$ws.Cells.Item(5, 4).Value = "'352.73"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +5.45%  "
$ws.Cells.Item(6, 5).Value = "  +0.43%  "
$ws.Cells.Item(7, 4).Value = "'0.5284"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +1.14%  "
$ws.Cells.Item(8, 5).Value = "  +1.14%  "
$ws.Cells.Item(9, 4).Value = "'54.31"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +0.99%  "
$ws.Cells.Item(10, 4).Value = "'0.09193"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +3.01%  "
$ws.Cells.Item(11, 4).Value = "'1.185"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +0.17%  "
$ws.Cells.Item(12, 4).Value = "'24.97"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +3.40%  "
$ws.Cells.Item(13, 4).Value = "2.147.36"
$ws.Cells.Item(13, 5).Value = "  +1.67%  "
$ws.Cells.Item(14, 4).Value = "'6.915"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +1.45%  "
$ws.Cells.Item(15, 4).Value = "'8.174"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +1.98%  "
$ws.Cells.Item(16, 4).Value = "'102.52"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +5.98%  "
$ws.Cells.Item(17, 4).Value = "'0.00001181"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +3.47%  "
$ws.Cells.Item(18, 4).Value = "'1.012"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +0.40%  "
$ws.Cells.Item(19, 4).Value = "'0.06732"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +1.27%  "
$ws.Cells.Item(20, 4).Value = "'19.62"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +2.20%  "
$ws.Cells.Item(21, 5).Value = "  +0.29%  "
$ws.Cells.Item(22, 4).Value = "'6.379"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +0.94%  "
$ws.Cells.Item(23, 4).Value = "30.815.13"
$ws.Cells.Item(23, 5).Value = "  +0.82%  "
$ws.Cells.Item(24, 4).Value = "'12.89"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +3.46%  "
$ws.Cells.Item(25, 4).Value = "'2.386"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +1.44%  "
$ws.Cells.Item(26, 4).Value = "2.380.60"
$ws.Cells.Item(26, 5).Value = "  +1.05%  "
$ws.Cells.Item(27, 4).Value = "'22.61"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +1.96%  "
$ws.Cells.Item(28, 4).Value = "'2.632"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +4.25%  "
$ws.Cells.Item(29, 4).Value = "'165.23"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +1.60%  "
$ws.Cells.Item(30, 4).Value = "'136.96"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +2.62%  "
$ws.Cells.Item(31, 4).Value = "'1.222"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +1.10%  "
$ws.Cells.Item(32, 5).Value = "  +1.30%  "
$ws.Cells.Item(33, 4).Value = "'1.671"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +2.13%  "
$ws.Cells.Item(34, 4).Value = "'6.424"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +0.04%  "
$ws.Cells.Item(35, 5).Value = "  +1.81%  "
$ws.Cells.Item(36, 4).Value = "'6.145"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +6.05%  "
$ws.Cells.Item(37, 4).Value = "'10.46"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +0.20%  "
$ws.Cells.Item(38, 4).Value = "'0.02662"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +3.04%  "
$ws.Cells.Item(39, 4).Value = "'0.06932"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +1.17%  "
$ws.Cells.Item(40, 4).Value = "'0.2338"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +1.90%  "
$ws.Cells.Item(41, 5).Value = "  +0.00%  "
$ws.Cells.Item(42, 4).Value = "'0.7002"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +1.95%  "
$ws.Cells.Item(43, 4).Value = "'1.274"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +2.13%  "
$ws.Cells.Item(44, 4).Value = "'14.83"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +5.94%  "
$ws.Cells.Item(45, 4).Value = "'2.368"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +2.14%  "
$ws.Cells.Item(46, 4).Value = "'0.6491"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +1.93%  "
$ws.Cells.Item(47, 4).Value = "'0.00000000371"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +6.68%  "
$ws.Cells.Item(48, 4).Value = "'3.759"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +2.56%  "
$ws.Cells.Item(49, 4).Value = "'1.261"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +0.96%  "
$ws.Cells.Item(50, 4).Value = "'83.45"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +0.33%  "
$ws.Cells.Item(51, 4).Value = "'0.07327"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +2.57%  "
